$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 4 inputs
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Home"
$ws.Range("Q4").Value = "AgileExtract2.csv"
$ws.Range("AF4").Value = "Home"
$ws.Range("AL4").Value = 1

# Update the selection to match the diff (S20)
$ws.Range("S20").Select()
